# Applies the "Penalty Reward System" (unfinished) edit:
#  - On "Weekly Quantity": remove the weekly rows for 2023-07-09, 2023-07-23,
#    2023-07-30, 2024-03-10, 2024-03-17 and 2024-03-24 (rows 15, 17, 18, 30,
#    31, 32), shifting remaining rows up.
#  - On "Monthly Trend": change the July 2023 quantity (row 6, col B) from
#    350 to 84, and remove the March 2024 row (original row 13), shifting
#    remaining rows up.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Delete rows from bottom to top so row indices of earlier rows stay valid.
$ws1.Rows.Item(32).Delete()
$ws1.Rows.Item(31).Delete()
$ws1.Rows.Item(30).Delete()
$ws1.Rows.Item(18).Delete()
$ws1.Rows.Item(17).Delete()
$ws1.Rows.Item(15).Delete()

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update the requested quantity for the 2023-07-31 month row.
$ws2.Range("B6").Value = 84

# Delete the 2024-03-31 month row.
$ws2.Rows.Item(13).Delete()
